# class_progression_Champion.xlsx update
#
# The "Advancement" columns that had no feat/ability listed for a given
# level were left as blank (but styled) cells. This script fills every
# one of those trailing blanks with the literal text "none" so the XML
# map / table backing the sheet always has a value to bind to.
#
# (Previously blank ranges, one contiguous block per row from the first
# empty Advancement column through column K.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blankRanges = @(
    "D3:K3",
    "E4:K4",
    "D5:K5",
    "F6:K6",
    "D7:K7",
    "F8:K8",
    "D9:K9",
    "H10:K10",
    "E11:K11",
    "G12:K12",
    "D13:K13",
    "F14:K14",
    "D15:K15",
    "F16:K16",
    "D17:K17",
    "F18:K18",
    "D19:K19",
    "E20:K20",
    "E21:K21"
)

foreach ($rangeAddress in $blankRanges) {
    $ws.Range($rangeAddress).Value = "none"
}

# Move the active selection (matches the author's last on-screen state).
$ws.Range("M3").Select()
